$wb = $excel.ActiveWorkbook

# New "Hora" timestamps for each forecast sheet (rows 2-25, column A).
# Each sheet holds 24 hourly timestamps counting back from a different
# reference hour: 24_HRS = now, 1d_bef = 1 day back, 7d_bef = 7 days back.
$sheetValues = @{
    "24_HRS" = @("2023-05-12 07:30", "2023-05-12 06:30", "2023-05-12 05:30", "2023-05-12 04:30", "2023-05-12 03:30", "2023-05-12 02:30", "2023-05-12 01:30", "2023-05-12 00:30", "2023-05-11 23:30", "2023-05-11 22:30", "2023-05-11 21:30", "2023-05-11 20:30", "2023-05-11 19:30", "2023-05-11 18:30", "2023-05-11 17:30", "2023-05-11 16:30", "2023-05-11 15:30", "2023-05-11 14:30", "2023-05-11 13:30", "2023-05-11 12:30", "2023-05-11 11:30", "2023-05-11 10:30", "2023-05-11 09:30", "2023-05-11 08:30")
    "1d_bef" = @("2023-05-11 08:30", "2023-05-11 07:30", "2023-05-11 06:30", "2023-05-11 05:30", "2023-05-11 04:30", "2023-05-11 03:30", "2023-05-11 02:30", "2023-05-11 01:30", "2023-05-11 00:30", "2023-05-10 23:30", "2023-05-10 22:30", "2023-05-10 21:30", "2023-05-10 20:30", "2023-05-10 19:30", "2023-05-10 18:30", "2023-05-10 17:30", "2023-05-10 16:30", "2023-05-10 15:30", "2023-05-10 14:30", "2023-05-10 13:30", "2023-05-10 12:30", "2023-05-10 11:30", "2023-05-10 10:30", "2023-05-10 09:30")
    "7d_bef" = @("2023-05-05 08:30", "2023-05-05 07:30", "2023-05-05 06:30", "2023-05-05 05:30", "2023-05-05 04:30", "2023-05-05 03:30", "2023-05-05 02:30", "2023-05-05 01:30", "2023-05-05 00:30", "2023-05-04 23:30", "2023-05-04 22:30", "2023-05-04 21:30", "2023-05-04 20:30", "2023-05-04 19:30", "2023-05-04 18:30", "2023-05-04 17:30", "2023-05-04 16:30", "2023-05-04 15:30", "2023-05-04 14:30", "2023-05-04 13:30", "2023-05-04 12:30", "2023-05-04 11:30", "2023-05-04 10:30", "2023-05-04 09:30")
}

foreach ($sheetName in $sheetValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $sheetValues[$sheetName]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $values[$i]
    }
}
